$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename "roundabout road 1 lane" -> "junction road 1 lane"
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("roundabout road 1 lane")
$ws2.Name = "junction road 1 lane"

# ------------------------------------------------------------------
# 2. Insert a new column D ("Vehicle count in Plus Junction in last
#    time stamp"), shifting the old D:I speed columns to E:J
# ------------------------------------------------------------------
$ws2.Columns.Item(4).Insert() | Out-Null
$ws2.Range("D1").Value = "Vehicle count in Plus Junction in last time stamp"

# ------------------------------------------------------------------
# 3. Remove the old rows 4 & 5 (only 2 data rows remain now)
# ------------------------------------------------------------------
$ws2.Rows.Item(4).Resize(2).Delete() | Out-Null

# ------------------------------------------------------------------
# 4. Overwrite the remaining two data rows with the new values
# ------------------------------------------------------------------
$ws2.Range("A2").Value = "2024-09-05 22:00:26"
$ws2.Range("B2").Value = 41.5556004701413
$ws2.Range("C2").Value = 10
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 45.46784310177639
$ws2.Range("F2").Value = 49.77413081422966
$ws2.Range("G2").Value = 41.12680021475725
$ws2.Range("H2").Value = 41.93699337804463
$ws2.Range("I2").Value = 50.40583294549832
$ws2.Range("J2").Value = 3.600646894576797

$ws2.Range("A3").Value = "2024-09-05 22:00:31"
$ws2.Range("B3").Value = 37.55672800173495
$ws2.Range("C3").Value = 17
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 45.46784310177639
$ws2.Range("F3").Value = 44.99773386020145
$ws2.Range("G3").Value = 34.86979835085188
$ws2.Range("H3").Value = 42.4491405187034
$ws2.Range("I3").Value = 37.58831652549469
$ws2.Range("J3").Value = 15.34775738813367

# ------------------------------------------------------------------
# 5. Add the new sheet "straight road 2 lanes" after the last sheet
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "straight road 2 lanes"

# Copy the header formatting (bold font, borders, center/top align)
# used on the other sheets so the new header row matches it.
$wb.Worksheets.Item("Accidents").Range("A1").Copy() | Out-Null
$ws3.Range("A1:E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws3.Range("A1").Value = "Time Stamp"
$ws3.Range("B1").Value = "Average Speed"
$ws3.Range("C1").Value = "Density"
$ws3.Range("D1").Value = "Avg. Speed (Road 0, Direction 0)"
$ws3.Range("E1").Value = "Avg. Speed (Road 0, Direction 1)"

$ws3.Range("A2").Value = "2024-09-05 22:00:41"
$ws3.Range("B2").Value = 72.17388803611196
$ws3.Range("C2").Value = 12
$ws3.Range("D2").Value = 75.82425041078419
$ws3.Range("E2").Value = 68.52352566143973

$ws3.Range("A3").Value = "2024-09-05 22:00:46"
$ws3.Range("B3").Value = 73.14244518124306
$ws3.Range("C3").Value = 19
$ws3.Range("D3").Value = 76.00193087164676
$ws3.Range("E3").Value = 69.96523885857228

Write-Host "edit complete"
